$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns J (Bought), K (Used), L (Available) ---
# ("Notes" header text is added after the note cells below, to line up
#  with the shared-string order of the authored workbook.)
$ws.Range("J1").Value = "Bought"
$ws.Range("K1").Value = "Used"
$ws.Range("L1").Value = "Available"

# Row 2
$ws.Range("J2").Value = 1
$ws.Range("K2").Formula = "=F2"
$ws.Range("L2").Formula = "=MAX(0,J2-K2)"

# Row 3
$ws.Range("J3").Value = 10
$ws.Range("K3").Formula = "=F3"
$ws.Range("L3").Formula = "=MAX(0,J3-K3)"

# Row 4
$ws.Range("J4").Value = 10
$ws.Range("K4").Formula = "=F4"
$ws.Range("L4").Formula = "=MAX(0,J4-K4)"

# Row 5
$ws.Range("J5").Value = 10
$ws.Range("K5").Formula = "=F5"
$ws.Range("L5").Formula = "=MAX(0,J5-K5)"

# Row 6
$ws.Range("J6").Value = 9
$ws.Range("K6").Formula = "=F6"
$ws.Range("L6").Formula = "=MAX(0,J6-K6)"

# Row 7
$ws.Range("J7").Value = 10
$ws.Range("K7").Formula = "=F7"
$ws.Range("L7").Formula = "=MAX(0,J7-K7)"

# Row 8
$ws.Range("J8").Value = 10
$ws.Range("K8").Formula = "=F8"
$ws.Range("L8").Formula = "=MAX(0,J8-K8)"

# Row 9
$ws.Range("J9").Value = 4
$ws.Range("K9").Formula = "=F9"
$ws.Range("L9").Formula = "=MAX(0,J9-K9)"

# Row 10
$ws.Range("J10").Value = 5
$ws.Range("K10").Formula = "=F10"
$ws.Range("L10").Formula = "=MAX(0,J10-K10)"
$ws.Range("M10").Value = "2 more on breakout boards"

# Row 11
$ws.Range("J11").Value = 2
$ws.Range("K11").Formula = "=F11"
$ws.Range("L11").Formula = "=MAX(0,J11-K11)"

# Row 12
$ws.Range("J12").Value = 1
$ws.Range("K12").Formula = "=F12"
$ws.Range("L12").Formula = "=MAX(0,J12-K12)"

# Row 13
$ws.Range("J13").Value = 2
$ws.Range("K13").Formula = "=F13"
$ws.Range("L13").Formula = "=MAX(0,J13-K13)"

# Row 14
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = 5
$ws.Range("L14").Formula = "=MAX(0,J14-K14)"

# Row 15
$ws.Range("J15").Value = 10
$ws.Range("K15").Formula = "=F15"
$ws.Range("L15").Formula = "=MAX(0,J15-K15)"

# Row 16
$ws.Range("J16").Value = 11
$ws.Range("K16").Formula = "=F16"
$ws.Range("L16").Formula = "=MAX(0,J16-K16)"

# Row 17
$ws.Range("J17").Value = 16
$ws.Range("K17").Formula = "=F17"
$ws.Range("L17").Formula = "=MAX(0,J17-K17)"

# Row 18
$ws.Range("J18").Value = 10
$ws.Range("K18").Formula = "=F18"
$ws.Range("L18").Formula = "=MAX(0,J18-K18)"

# Row 19
$ws.Range("J19").Value = 10
$ws.Range("K19").Formula = "=F19"
$ws.Range("L19").Formula = "=MAX(0,J19-K19)"

# Row 20
$ws.Range("J20").Value = 10
$ws.Range("K20").Value = 3
$ws.Range("L20").Formula = "=MAX(0,J20-K20)"

# Row 21
$ws.Range("J21").Value = 2
$ws.Range("K21").Formula = "=F21"
$ws.Range("L21").Formula = "=MAX(0,J21-K21)"

# Row 22
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 0
$ws.Range("L22").Formula = "=MAX(0,J22-K22)"
$ws.Range("M22").Value = "Used 100K instead"

# --- Notes header (added last so it lands after the note strings above) ---
$ws.Range("M1").Value = "Notes"

# --- Match formatting of existing table (thin border box) on the new J:M columns ---
$ws.Range("J1:M22").Borders.LineStyle = 1

# --- Column M width (matches author's manual resize / best-fit) ---
$ws.Columns("M").ColumnWidth = 25.140625

# --- Move active selection like in the authored workbook ---
$ws.Range("J16").Select()
